# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Thu May  9 06:30:32 UTC 2024 with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS like a plain decimal number (e.g. "1.00", "7.02")
# as literal TEXT, matching the workbook's original inlineStr/shared-string cells.
# Forcing NumberFormat="@" prevents Excel from auto-coercing the string to a number
# (which would silently drop the trailing zero, e.g. "1.00" -> 1). ClearFormats()
# afterwards restores the cell's original (default) style so only the value changes.
function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}

# --- Plain text updates (Coin/Link renames, percentage volumes, and price strings ---
# --- that already fail numeric parsing, e.g. thousand-dot formatted "61.714.72") ----
$ws.Range("D2").Value = "61.714.72"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "3.004.89"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "3.002.92"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  +6.41%  "
$ws.Range("E12").Value = "  +4.89%  "
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "3.500.47"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "61.617.76"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "3.003.68"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("E24").Value = "  +2.73%  "
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  +9.08%  "
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E30").Value = "  +2.99%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("D35").Value = "0.0₃0844"
$ws.Range("E35").Value = "  +8.23%  "
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  +11.15%  "
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").Value = "2.712.37"
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("E51").Value = "  +3.00%  "

# --- Price updates that look like plain decimals: forced to stay TEXT -------------
Set-TextValue "D4" "1.00"
Set-TextValue "D5" "599.71"
Set-TextValue "D6" "145.03"
Set-TextValue "D8" "0.522"
Set-TextValue "D10" "0.148"
Set-TextValue "D11" "6.04"
Set-TextValue "D14" "34.52"
Set-TextValue "D17" "7.02"
Set-TextValue "D20" "452.44"
Set-TextValue "D21" "14.06"
Set-TextValue "D22" "0.688"
Set-TextValue "D23" "7.37"
Set-TextValue "D24" "81.80"
Set-TextValue "D25" "2.25"
Set-TextValue "D27" "12.03"
Set-TextValue "D30" "7.27"
Set-TextValue "D31" "1.00"
Set-TextValue "D33" "27.51"
Set-TextValue "D36" "1.02"
Set-TextValue "D37" "5.79"
Set-TextValue "D40" "50.40"
Set-TextValue "D41" "0.124"
Set-TextValue "D42" "2.93"
Set-TextValue "D43" "403.19"
Set-TextValue "D44" "39.69"
Set-TextValue "D45" "0.272"
Set-TextValue "D46" "0.0355"
Set-TextValue "D48" "131.89"

